# Auto-generated Excel COM-interop script
# Applies the "pu" sheet data-model refresh + related workbook metadata edits

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pu")

# --- Rewrite the data block (A2:D58) in one shot -------------------------
# Clear everything below the header first so stale cells beyond the new
# data extent (columns/rows that no longer hold data) are removed.
$ws.Range("A2:H300").ClearContents()

# Re-stamp column A (the "id" column) with its existing fill/number format
# before writing values, so every id cell keeps the s="15" style used
# throughout the column.
$ws.Range("A2").Copy()
$ws.Range("A2:A58").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$data = New-Object 'object[,]' 57,4
$data[0,0] = 1; $data[0,1] = 60.410000000000004; $data[0,2] = -28.66023; $data[0,3] = 21.13877
$data[1,0] = 2; $data[1,1] = 48.3; $data[1,2] = -28.66422; $data[1,3] = 21.13672
$data[2,0] = 3; $data[2,1] = 73; $data[2,2] = -28.781; $data[2,3] = 20.781
$data[3,0] = 4; $data[3,1] = 52.019999999999996; $data[3,2] = -28.43161; $data[3,3] = 20.65984
$data[4,0] = 5; $data[4,1] = 42.96; $data[4,2] = -28.65328; $data[4,3] = 21.15693
$data[5,0] = 6; $data[5,1] = 32.339999999999996; $data[5,2] = -28.65777; $data[5,3] = 21.125
$data[6,0] = 7; $data[6,1] = 61.440000000000005; $data[6,2] = -28.4611; $data[6,3] = 20.65961
$data[7,0] = 8; $data[7,1] = 13.35; $data[7,2] = -28.80554; $data[7,3] = 20.65364
$data[8,0] = 9; $data[8,1] = 100.27; $data[8,2] = -28.498; $data[8,3] = 20.144
$data[9,0] = 10; $data[9,1] = 37.73; $data[9,2] = -28.49469; $data[9,3] = 20.14799
$data[10,0] = 11; $data[10,1] = 29.47; $data[10,2] = -28.43078; $data[10,3] = 20.14013
$data[11,0] = 12; $data[11,1] = 100; $data[11,2] = -28.455; $data[11,3] = 20.044
$data[12,0] = 14; $data[12,1] = $null; $data[12,2] = -28.6542037; $data[12,3] = 19.5236904
$data[13,0] = 15; $data[13,1] = 42.585; $data[13,2] = -28.6542037; $data[13,3] = 19.5236904
$data[14,0] = 16; $data[14,1] = 49.273; $data[14,2] = -28.6542037; $data[14,3] = 19.5236904
$data[15,0] = 17; $data[15,1] = $null; $data[15,2] = -28.6542037; $data[15,3] = 19.5236904
$data[16,0] = 18; $data[16,1] = 46.85; $data[16,2] = -28.6542037; $data[16,3] = 19.5236904
$data[17,0] = 19; $data[17,1] = $null; $data[17,2] = -28.6542037; $data[17,3] = 19.5236904
$data[18,0] = 20; $data[18,1] = $null; $data[18,2] = -28.6542037; $data[18,3] = 19.5236904
$data[19,0] = 21; $data[19,1] = 45.11; $data[19,2] = -28.6542037; $data[19,3] = 19.5236904
$data[20,0] = 22; $data[20,1] = 47.78; $data[20,2] = -28.6542037; $data[20,3] = 21.11212
$data[21,0] = 23; $data[21,1] = 43.750000000000014; $data[21,2] = -28.6542037; $data[21,3] = 19.5236904
$data[22,0] = 24; $data[22,1] = 40.38; $data[22,2] = -28.9583855; $data[22,3] = 19.0004562
$data[23,0] = 25; $data[23,1] = 44.849999999999994; $data[23,2] = -28.9583855; $data[23,3] = 19.0004562
$data[24,0] = 26; $data[24,1] = 15.04; $data[24,2] = -28.9583855; $data[24,3] = 19.0004562
$data[25,0] = 27; $data[25,1] = $null; $data[25,2] = -28.6542037; $data[25,3] = 19.5236904
$data[26,0] = 33; $data[26,1] = $null; $data[26,2] = -34.076157; $data[26,3] = 18.892171
$data[27,0] = 43; $data[27,1] = 51.480000000000004; $data[27,2] = -28.65393; $data[27,3] = 21.11125
$data[28,0] = 44; $data[28,1] = 100.33999999999996; $data[28,2] = -33.444967; $data[28,3] = 19.629232
$data[29,0] = 45; $data[29,1] = 61.84; $data[29,2] = -28.66343; $data[29,3] = 21.15783
$data[30,0] = 46; $data[30,1] = 55.41; $data[30,2] = -28.66813; $data[30,3] = 21.15195
$data[31,0] = 47; $data[31,1] = 61.79999999999999; $data[31,2] = -28.65414; $data[31,3] = 21.15917
$data[32,0] = 48; $data[32,1] = 66.14; $data[32,2] = -33.793974; $data[32,3] = 19.828017
$data[33,0] = 73; $data[33,1] = 8; $data[33,2] = -28.65372; $data[33,3] = 21.11456
$data[34,0] = 74; $data[34,1] = 71.78999999999999; $data[34,2] = -32.901697; $data[34,3] = 18.746009
$data[35,0] = 75; $data[35,1] = 40.620000000000005; $data[35,2] = -33.793974; $data[35,3] = 19.828017
$data[36,0] = 76; $data[36,1] = 60.61; $data[36,2] = -32.901697; $data[36,3] = 18.746009
$data[37,0] = 93; $data[37,1] = 114.86; $data[37,2] = -25.093; $data[37,3] = 29.387
$data[38,0] = 77; $data[38,1] = $null; $data[38,2] = -28.6542037; $data[38,3] = 19.5236904
$data[39,0] = 78; $data[39,1] = $null; $data[39,2] = -28.6542037; $data[39,3] = 19.5236904
$data[40,0] = 79; $data[40,1] = $null; $data[40,2] = -28.6542037; $data[40,3] = 19.5236904
$data[41,0] = 80; $data[41,1] = $null; $data[41,2] = -28.6542037; $data[41,3] = 19.5236904
$data[42,0] = 81; $data[42,1] = $null; $data[42,2] = -28.6542037; $data[42,3] = 19.5236904
$data[43,0] = 82; $data[43,1] = $null; $data[43,2] = -28.6542037; $data[43,3] = 19.5236904
$data[44,0] = 83; $data[44,1] = $null; $data[44,2] = -28.6542037; $data[44,3] = 19.5236904
$data[45,0] = 84; $data[45,1] = $null; $data[45,2] = -28.6542037; $data[45,3] = 19.5236904
$data[46,0] = 85; $data[46,1] = $null; $data[46,2] = -28.6542037; $data[46,3] = 19.5236904
$data[47,0] = 68; $data[47,1] = $null; $data[47,2] = -28.647; $data[47,3] = 19.515
$data[48,0] = 87; $data[48,1] = $null; $data[48,2] = -28.647; $data[48,3] = 19.515
$data[49,0] = 88; $data[49,1] = $null; $data[49,2] = -28.647; $data[49,3] = 19.515
$data[50,0] = 92; $data[50,1] = $null; $data[50,2] = -28.66023; $data[50,3] = 21.13877
$data[51,0] = 96; $data[51,1] = $null; $data[51,2] = -25.093; $data[51,3] = 29.387
$data[52,0] = 55; $data[52,1] = $null; $data[52,2] = -28.795; $data[52,3] = 20.64
$data[53,0] = 51; $data[53,1] = $null; $data[53,2] = -28.795; $data[53,3] = 20.64
$data[54,0] = 52; $data[54,1] = $null; $data[54,2] = -28.795; $data[54,3] = 20.64
$data[55,0] = 66; $data[55,1] = $null; $data[55,2] = -28.795; $data[55,3] = 20.64
$data[56,0] = 69; $data[56,1] = $null; $data[56,2] = -28.795; $data[56,3] = 20.64

$ws.Range("A2:D58").Value2 = $data

# --- New helper cells H37:H38 (litres/min -> per-hour conversion) --------
$ws.Range("H37").Value2 = 0.05
$ws.Range("H38").Formula = "=H37*60"

# --- New row 39 (id 93) gets right-aligned number formatting on column A -
$ws.Range("A39").HorizontalAlignment = -4152  # xlRight

# --- Row 52 C:D picked up the "wrap text / vertical center" style that was
# already defined in the workbook (index 1) when the values were pasted in -
$ws.Range("C52:D52").VerticalAlignment = -4108  # xlCenter
$ws.Range("C52:D52").WrapText = $true

# --- Sheet view: scrolled down to show the newly appended rows, selection
# parked on the helper cell that was last edited ---------------------------
$win = $wb.Windows.Item(1)
$ws.Activate()
$ws.Range("A22").Select()
$win.ScrollRow = 22
$ws.Range("H39").Select()

# --- AutoFilter now spans the larger A1:D38 footprint (filter header row +
# the underlying ids, rather than covering every appended lookup row) ------
$ws.AutoFilterMode = $false
$ws.Range("A1:D38").AutoFilter()

# --- Printable page setup (paper size / orientation) now explicit ---------
$ws.PageSetup.PaperSize = 9   # xlPaperA4
$ws.PageSetup.Orientation = 1 # xlPortrait

# --- Workbook-level bookkeeping: the "pu" tab was recreated upstream in the
# source database-refresh pipeline so defined-name ranges referencing it need
# to track the new data extent -----------------------------------------------
try {
    $wb.Names.Item("_xlnm._FilterDatabase").RefersToR1C1 = "=pu!R1C1:R38C4"
} catch { }

